$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: new article entry (Ser 33, published 2020-02-02)
# Set string-valued cells in the same order the original workbook's
# shared-string table gained them (Ayats, Tags, then Content) so the
# shared string indices line up: 90=Ayats, 91=Tags, 92=Content.
$ws.Cells.Item(34, 3).Value = "Surah Al Imran, 65 - 83"
$ws.Cells.Item(34, 6).Value = "Colorful vibes, Nuclear War, Life Expectancy, Social Center, Speed of Time, Decision Making"
$ws.Cells.Item(34, 4).Value = @"
h1: How long does it take to live?
p.note: This topic came to my mind a moment ago standing in kitchen making coffee for wife. Ok, with coffee..
p: Guys, last year on 27th February 2019, it was a bold moment for all of us. We downed an Indian aircraft. A day prior, we were really annoyed for letting Indian jets fly back out with out targeting it. It was an attack on our sovereignty. 26th February 2019, we were angry. 27th February 2019, the next day we were proud. Today is 2 February, we are sleeping fine and forgetful, how close we went to a nuclear war. It has been 11 months since last year’s major escalation. 
p: These past 11 months, there had been numerous fights, a number of losses and so many failures. We even lost people very dear to us. After all these downs, we are still sitting here in our cozy beds, not much affected. 
p: Simply, we lived a year of troubles without even feeling a breeze of it passing. The brink of nuclear war failed to wake us up. Today Coronavirus is spreading like a fire storm, but we are sure we can not die through this plague. Deep in heart, we know it is in China not in my city.
h3: Life expectancy?
p: A good healthy Pakistani lives around 70 to 80 years if he does not die due to a heart attack from extra spicy desi dishes. A good athlete lives 100 years, still with all chances taken and no natural calamity occurring. Most die of earth quakes and car accidents. 
p: 9 Million people die each year due to cardiac arrests. 1.3 Million people die each year due to road accidents. <a href= target=_blanl>Link here.</a> 99.999 % die each year over the age of 100 years old. I am 32 right now. If I am lucky I will die in 60 years. 
p: Earth came into being about 50 Million years back. So many generations have lived since earth’s creation, affected by worst storms. Just like the famous movie, 2012, Nooh (as)’s people have seen the entire planet flooding with water. They sit in a ship and survive. Some generations were freed of Pharoah and some were saved from a complete Nuclear war on 27 February 2020. 
h3: Why are we sitting dumb folded?
p: May be we need to put a reminder beside our bed or write in bold big letters somewhere prominent in our homes about this time leaving us behind. Else this year is also going to pass and we will still be sitting on cross roads looking at sky wondering when the next buzz will wake us up.
p: We are moving quite fast each day. Last year it was escalation, we were standing on the brink of a nuclear war. Then we talked of North Korea fighting US. Now we are shouting it is Iran vs Arab vs US. Time is simply passing by us. It is really ruthless, doesn’t care if we are careful in handling it or not. It is hitting us with failures, losses and hurdles that we wake up, but we are resorting to being more safer and less riskier. We are not kicking our adrenaline, we are not becoming adventurous and we are not worried about another year spent in loss. 
p: Guys, the best way is the Quran way. It became all adventurous and still it is focusing on us to liberate ourselves of our own fears. Willingly or un-willingly we are going to meet our Creator in some moments, it is upto us how we prepare ourselves.
quote: So is it other than the religion of Allah they desire, while to Him have submitted [all] those within the heavens and earth, willingly or by compulsion, and to Him they will be returned? <br> - Surah Al Imran verse 83
p.b-left: Do we stick to same old rotten legacies laid out by our ancestors or we go out find the new bright truth?
p.b-left: Do we even plan to take some falls this year?
p.b-left: Do we want to portray ourselves an example for our kids?
p.b-left: Do we really care about the time rushing through us?
p.b-left: How are we going to stand and when? What is going to come to make us live a bolder and a better life? 
p: I am talking all macho right now, I know I lived an entirety in darkness. This past month is a roller coaster of red bold high events all under the shadow of Quran. Each day is zillion times more productive than yesterday. 
p: I really want we all taste the colors of our mosques. It is library today, laboratory tomorrow and a place to play with our kids the rest of the year. We all have mosques in our towns, that can serve as an amazing co-working space, with a projector screen and sermons being delivered using laptops, nicely built powerpoint presentations. People are discussing ideas in already built calm mosques serving as the center of our light solving our basic social issues. 
p: Mosques have the potential to repair the damaged part of our society. We need more engineers, doctors and scientists to lit them up, paint them with good colors, install large video walls and nice spot lights that attracts our kids.
h3: Conclusion
p: It really takes a moment to change our entire life. 1 single moment made me a programmer where I decided to purchase a 11 Dollar Web Development Course on Udemy in 2016. 1 single moment made me smoker-less where I decided I will not move away from my mother for 3 days. 1 single moment made me rich where I decided I will no more bound my self by illogical limits. It only takes a moment. The blink of an eye is all it takes to live. 
p.note: Guys, if you like this project. Please follow this project's page on twitter. <a href="https://twitter.com/zakatlists">Click here to go to the twitter page</a>.
"@
$ws.Cells.Item(34, 5).Value = "Qasim Ali"

# Ser number
$ws.Cells.Item(34, 1).Value = 33

# Date (2020-02-02), matching the date-formatted style already used by column B
$ws.Cells.Item(34, 2).Value = 43863
$ws.Cells.Item(34, 2).NumberFormat = "d-mmm-yy"

# Row height: long wrapped content fills the row to Excel's max row height,
# same as every other fully-populated row in this sheet.
$ws.Rows.Item(34).RowHeight = 409.6

# Move the active selection to the newly filled row, like the author did.
$ws.Range("D34").Select() | Out-Null
